# Weekly price-sheet update: insert the newest "Achicoria" (Vega Central
# Mapocho de Santiago) price observation as a new row right after the
# existing header block (row 14), pushing all prior data rows down by one.
# This mirrors how the daily/weekly logic keeps the most-recent date in the
# earliest data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 15; everything currently at row 15 onward
# (through row 53) shifts down to rows 16-54.
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new weekly observation.
$ws.Cells.Item(15, 1).Value  = 9
$ws.Cells.Item(15, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(15, 3).Value  = "Metropolitana"
$ws.Cells.Item(15, 4).Value  = 45030
$ws.Cells.Item(15, 5).Value  = 13
$ws.Cells.Item(15, 6).Value  = 100112010
$ws.Cells.Item(15, 7).Value  = "Achicoria"
$ws.Cells.Item(15, 8).Value  = "Sin especificar"
$ws.Cells.Item(15, 9).Value  = "Primera"
$ws.Cells.Item(15, 10).Value = 90
$ws.Cells.Item(15, 11).Value = 7000
$ws.Cells.Item(15, 12).Value = 7000
$ws.Cells.Item(15, 13).Value = 7000
$ws.Cells.Item(15, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(15, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(15, 16).Value = 438
$ws.Cells.Item(15, 17).Value = 16
$ws.Cells.Item(15, 18).Value = "Hortaliza"
